# Apply edits: course group and custom group added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - replace names (course group)
$ws.Range("A4").Value = "John"
$ws.Range("B4").Value = "Henry"

# Row 5 - replace names (custom group) and make C5 a numeric ID instead of text
$ws.Range("A5").Value = "Shirish"
$ws.Range("B5").Value = "Dangol"
$ws.Range("C5").Value = 1234567

# Update the active selection to match the saved view state
$ws.Range("H12").Select()
